$d = $word.ActiveDocument

# The logo pictures embedded in the headers/footers were exported from
# PowerPoint with mismatched internal "name" labels (the docPr / cNvPr
# `name` attribute is just a cosmetic label Word shows in the Selection
# Pane - it does NOT control which media part is rendered, that's the
# blip r:embed relationship, which is untouched here). This renames the
# labels so they line up with the image's actual file extension:
#   word/footer1.xml, word/footer2.xml : image1.png -> image2.png
#   word/header1.xml, word/header2.xml : image2.jpg -> image1.jpg
#
# InlineShape exposes no working Name setter in the object model, so we
# go straight at the package XML via Document.WordOpenXML (get the flat
# OPC XML, patch the four <pkg:part> segments precisely, write it back).

$xml = $d.WordOpenXML

function Get-PartBounds($text, $partName) {
    $startTag = '<pkg:part pkg:name="' + $partName + '"'
    $start = $text.IndexOf($startTag)
    if ($start -lt 0) {
        throw "part not found: $partName"
    }
    $closeTag = "</pkg:part>"
    $end = $text.IndexOf($closeTag, $start) + $closeTag.Length
    return @($start, $end)
}

function Rename-ImageInPart($text, $partName, $oldName, $newName) {
    $bounds = Get-PartBounds $text $partName
    $start = $bounds[0]
    $end = $bounds[1]
    $segment = $text.Substring($start, $end - $start)

    $patched = $segment.Replace('name="' + $oldName + '"', 'name="' + $newName + '"')

    return $text.Substring(0, $start) + $patched + $text.Substring($end)
}

$xml = Rename-ImageInPart $xml "/word/footer1.xml" "image1.png" "image2.png"
$xml = Rename-ImageInPart $xml "/word/footer2.xml" "image1.png" "image2.png"
$xml = Rename-ImageInPart $xml "/word/header1.xml" "image2.jpg" "image1.jpg"
$xml = Rename-ImageInPart $xml "/word/header2.xml" "image2.jpg" "image1.jpg"

$d.WordOpenXML = $xml

Write-Output "renamed logo image labels in headers/footers"
